{"js": "// Add a new \"Body Text\" paragraph at the end of the document that discusses\n// the follow-up short communication, including a hyperlink to it.\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// First sentence fragment, ending with the open quote before the linked title.\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Interestingly, the same group that published the initial data set recently published a short communication entitled '\",\n  \"After\"\n);\nnewParagraph.style = \"Body Text\";\n\n// The hyperlinked article title.\nconst linkRange = newParagraph.insertText(\n  \"A tool for calculating concentration ratios from large environmental data sets\",\n  \"End\"\n);\nlinkRange.hyperlink =\n  \"http://www.sciencedirect.com/science/article/pii/S0265931X1830081X\";\n\n// Remainder of the paragraph, immediately after the linked text (no space).\nnewParagraph.insertText(\n  \". This publication outlines a tool developed in MS Excel that predicts biological accumulation of mine contaminants, given their proximity to environmental samples with certain concentrations of both metals and radionuclides. We'd expect others to start analyzing this current data set in light of this newly developed tool, since these many environmental measures can now been linked to biological outcomes.\",\n  \"End\"\n);\n\nawait context.sync();\n", "ps1": "# Add a new \"Body Text\" paragraph at the end of the document that discusses\n# the follow-up short communication, including a hyperlink to it.\n$d = $word.ActiveDocument\n\n# Start a new paragraph right after the current last paragraph.\n$r = $d.Paragraphs.Last.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newPara.Style = \"Body Text\"\n$pRange = $newPara.Range\n\n# First sentence fragment, ending with the open quote before the linked title.\n$pRange.Collapse(0)\n$pRange.InsertAfter(\"Interestingly, the same group that published the initial data set recently published a short communication entitled '\")\n$pRange.Collapse(0)\n\n# The hyperlinked article title.\n$linkRange = $pRange.Duplicate\n$linkRange.InsertAfter(\"A tool for calculating concentration ratios from large environmental data sets\")\n$d.Hyperlinks.Add($linkRange, \"http://www.sciencedirect.com/science/article/pii/S0265931X1830081X\") | Out-Null\n\n# Remainder of the paragraph, immediately after the linked text (no space).\n$pRange.Collapse(0)\n$pRange.SetRange($linkRange.End, $linkRange.End)\n$pRange.InsertAfter(\". This publication outlines a tool developed in MS Excel that predicts biological accumulation of mine contaminants, given their proximity to environmental samples with certain concentrations of both metals and radionuclides. We'd expect others to start analyzing this current data set in light of this newly developed tool, since these many environmental measures can now been linked to biological outcomes.\")\n"}
